$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.127881588408715, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.449980674824537)
    3 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    4 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    5 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    6 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    7 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 3.645393585217082)
    8 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    9 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B: TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C: d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D: K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E: IP
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G: sum
}
